$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the style of existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the data for columns I (I0) and J (IF), rows 2-9
$i0 = @(8, 6, 9, 8, 6, 6, 6, 2)
$if = @(8, 8, 9, 9, 8, 6, 6, 2)

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 9).Value = $i0[$row - 2]
    $ws.Cells.Item($row, 10).Value = $if[$row - 2]
}
